$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the old rows 2-29 into the new condensed rows 2-7, then
# clear out what used to be rows 8-29 so the used range shrinks to A1:A7.

$ws.Range("A2").Value = "(""City's Blessing"", ['Card'])"
$ws.Range("A3").Value = "('Elemental', ['Token Creature — Elemental', '1/1'])"
$ws.Range("A4").Value = "('Golem', ['Token Artifact Creature — Golem', '4/4'])"
$ws.Range("A5").Value = "('Huatli, Radiant Champion Emblem', ['Emblem — Huatli', 'Whenever a creature enters the battlefield under your control, you may draw a card.'])"
$ws.Range("A6").Value = "('Rivals of Ixalan Checklist', ['Card', '(You can mark this card to represent a double-faced card in your library or hand.)', '☐ Hadana’s Climb {1}{G}{U}', '☐ Journey to Eternity {1}{B}{G}', '☐ Path of Mettle {R}{W}', '☐ Profane Procession {1}{W}{B}', '☐ Storm the Vault {2}{U}{R}', '☐ Azor’s Gateway {2}', '☐ Golden Guardian {4}'])"
$ws.Range("A7").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"

$ws.Range("A8:A29").ClearContents()
